$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5 (Category) ---
$ws.Range("C5").Value = 50
$ws.Range("D5").Formula = "=C5+10"
$ws.Range("E5").Formula = "=D5+10"
$ws.Range("F5").Formula = "=E5+10"
$ws.Range("G5").Formula = "=F5+10"

# --- Row 6 (User) ---
$ws.Range("C6").Value = 100
$ws.Range("D6").Formula = "=C6+10"
$ws.Range("E6").Formula = "=D6+20"
$ws.Range("F6").Formula = "=E6+20"
$ws.Range("G6").Formula = "=F6+20"

# --- Row 8 (UserOnImage) ---
$ws.Range("C8").Value = 500
$ws.Range("D8").Formula = "=C8+500"
$ws.Range("E8").Formula = "=D8+750"
$ws.Range("F8").Formula = "=E8+750"
$ws.Range("G8").Formula = "=F8+750"

# --- New cell D12 ---
$ws.Range("D12").Value = 5

# --- Selection moved ---
$ws.Range("H21").Select()
